$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($Table, $Row, $Col, $OldText, $NewText) {
    $cell = $Table.Cell($Row, $Col)
    $rng = $cell.Range
    $current = $rng.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $OldText) {
        Write-Host "WARNING: cell ($Row,$Col) expected old text but found a mismatch"
    }
    $rng.Text = $NewText
}

Set-CellValue $t 1 1 "79÷3=" "68÷6="
Set-CellValue $t 1 2 "41÷4=" "18÷7="
Set-CellValue $t 1 3 "74÷2=" "86÷3="
Set-CellValue $t 1 4 "81÷6=" "86÷9="
Set-CellValue $t 1 5 "73÷4=" "51÷3="
Set-CellValue $t 5 1 "70÷4=" "32÷8="
Set-CellValue $t 5 2 "11÷5=" "58÷6="
Set-CellValue $t 5 3 "81÷5=" "55÷2="
Set-CellValue $t 5 4 "86÷3=" "21÷5="
Set-CellValue $t 5 5 "64÷6=" "63÷6="
Set-CellValue $t 9 1 "12÷8=" "25÷2="
Set-CellValue $t 9 2 "13÷9=" "77÷3="
Set-CellValue $t 9 3 "96÷7=" "77÷3="
Set-CellValue $t 9 4 "76÷7=" "53÷4="
Set-CellValue $t 9 5 "86÷6=" "80÷7="
Set-CellValue $t 13 1 "38÷3=" "48÷3="
Set-CellValue $t 13 2 "52÷7=" "55÷5="
Set-CellValue $t 13 3 "27÷3=" "14÷9="
Set-CellValue $t 13 4 "33÷5=" "17÷4="
Set-CellValue $t 13 5 "49÷9=" "60÷3="
Set-CellValue $t 17 1 "79÷3=" "30÷5="
Set-CellValue $t 17 2 "75÷9=" "68÷6="
Set-CellValue $t 17 3 "85÷9=" "61÷2="
Set-CellValue $t 17 4 "62÷6=" "79÷4="
Set-CellValue $t 17 5 "21÷9=" "69÷4="
